$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the three image-link columns -------------------
$ws.Range("F1").Value = "Image 1"
$ws.Range("G1").Value = "Image 2"
$ws.Range("H1").Value = "Image 3"

# --- New image links for the "S45" location row (row 26) -----------------
$ws.Range("F26").Value = "https://i.postimg.cc/1XmqZGMQ/In_front_of_S45.jpg"
$ws.Range("G26").Value = "https://i.postimg.cc/d3syK26N/In_front_of_S45_9.jpg"
$ws.Range("H26").Value = "https://i.postimg.cc/PJXwkWyR/In_front_of_S45_7.jpg"

$ws.Hyperlinks.Add($ws.Range("F26"), "https://i.postimg.cc/1XmqZGMQ/In_front_of_S45.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G26"), "https://i.postimg.cc/d3syK26N/In_front_of_S45_9.jpg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("H26"), "https://i.postimg.cc/PJXwkWyR/In_front_of_S45_7.jpg") | Out-Null

# --- Widen the new image columns (manual resize by the author) -----------
$ws.Columns.Item(6).ColumnWidth = 43.529947916666664
$ws.Columns.Item(7).ColumnWidth = 44.346354166666664
$ws.Columns.Item(8).ColumnWidth = 45.166666666666664

# --- Freeze the header row and restore selection --------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.SplitRow = 1
$win.Split = $false
$win.FreezePanes = $true
$ws.Range("H27").Select()
$ws.Range("C1").Select()
